$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.598.54"
$ws.Range("E2").Value = "  -3.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.850.35"
$ws.Range("E3").Value = "  -3.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -1.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.85"
$ws.Range("E5").Value = "  +2.40%  "
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4662"
$ws.Range("E7").Value = "  -3.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3925"
$ws.Range("E8").Value = "  -3.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.59"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07918"
$ws.Range("E10").Value = "  -3.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9864"
$ws.Range("E11").Value = "  -2.39%  "
$ws.Range("E12").Value = "  -5.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.944.88"
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.850"
$ws.Range("E14").Value = "  -3.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.032"
$ws.Range("E15").Value = "  -3.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06895"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.80"
$ws.Range("E17").Value = "  -4.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001008"
$ws.Range("E19").Value = "  -3.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.12"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.600.34"
$ws.Range("E22").Value = "  -3.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.404"
$ws.Range("E23").Value = "  -4.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.28"
$ws.Range("E24").Value = "  -5.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.197.61"
$ws.Range("E25").Value = "  +3.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.130"
$ws.Range("E26").Value = "  -2.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.52"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.45"
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.131"
$ws.Range("E29").Value = "  -5.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.016"
$ws.Range("E30").Value = "  -3.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.54"
$ws.Range("E31").Value = "  -2.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9832"
$ws.Range("E32").Value = "  -3.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09431"
$ws.Range("E33").Value = "  -2.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.373"
$ws.Range("E34").Value = "  -4.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.487"
$ws.Range("E35").Value = "  -1.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.350"
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06162"
$ws.Range("E37").Value = "  -3.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02204"
$ws.Range("E38").Value = "  -4.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.161"
$ws.Range("E39").Value = "  -2.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5708"
$ws.Range("E40").Value = "  -4.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.621"
$ws.Range("E41").Value = "  -3.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.15"
$ws.Range("E42").Value = "  -5.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1796"
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.373"
$ws.Range("E44").Value = "  -4.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.251"
$ws.Range("E45").Value = "  -2.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5403"
$ws.Range("E46").Value = "  -3.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.80"
$ws.Range("E47").Value = "  -5.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.07158"
$ws.Range("E48").Value = "  -4.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.911"
$ws.Range("E49").Value = "  -2.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.22"
$ws.Range("E50").Value = "  -3.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "42.94"
$ws.Range("E51").Value = "  +2.26%  "
